# Adding options to fit sex ratios
# - Controls!B2 (n_sims): 10 -> 100
# - Controls!B6 (n_years): 31 -> 51
# - Growth_Param!A2 (k, Female): 0.15 -> 0.17
# - Growth_Param!B2 (k, Male): 0.1275 -> 0.2
# - Growth_Param!B3 (L_inf, Male): 80.75 -> 80
# - Recruitment_Mortality!B7 (sexRatio): 0.125 -> 0.15
# Also mirrors the author's final UI state: Controls!B2 selected, then
# Growth_Param activated with B4 selected (becomes the active tab).

$wb = $excel.ActiveWorkbook

# --- Controls sheet -------------------------------------------------
$wsControls = $wb.Worksheets.Item("Controls")
$wsControls.Activate()
$wsControls.Range("B2").Value = 100
$wsControls.Range("B6").Value = 51
$wsControls.Range("B2").Select()

# --- Growth_Param sheet ----------------------------------------------
$wsGrowth = $wb.Worksheets.Item("Growth_Param")
$wsGrowth.Activate()
$wsGrowth.Range("A2").Value = 0.17
$wsGrowth.Range("B2").Value = 0.2
$wsGrowth.Range("B3").Value = 80
$wsGrowth.Range("B4").Select()

# --- Recruitment_Mortality sheet -------------------------------------
$wsRecruit = $wb.Worksheets.Item("Recruitment_Mortality")
$wsRecruit.Range("B7").Value = 0.15

# Growth_Param is the sheet left active/selected in the saved workbook.
$wsGrowth.Activate()
$wsGrowth.Range("B4").Select()
